# fix: ensure elective courses are scheduled in same time slots for both sections A and B

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Section_A timetable updates
# ---------------------------------------------------------------------------
$wsA = $wb.Worksheets.Item("Section_A")

$wsA.Range("B2").Value = "CS263"
$wsA.Range("D2").Value = "Free"
$wsA.Range("E2").Value = "MA261"
$wsA.Range("F2").Value = "CS263"

$wsA.Range("B3").Value = "Free"
$wsA.Range("C3").Value = "Free"
$wsA.Range("E3").Value = "Free"
$wsA.Range("F3").Value = "Free"

$wsA.Range("C5").Value = "CS261"
$wsA.Range("D5").Value = "CS263"
$wsA.Range("F5").Value = "CS264 (Tutorial)"

$wsA.Range("B6").Value = "Free"
$wsA.Range("C6").Value = "MA261"
$wsA.Range("D6").Value = "CS261"
$wsA.Range("E6").Value = "CS261"
$wsA.Range("F6").Value = "Free"

$wsA.Range("B7").Value = "CS264"
$wsA.Range("C7").Value = "CS264"
$wsA.Range("E7").Value = "Free"
$wsA.Range("F7").Value = "CS264"

# ---------------------------------------------------------------------------
# Section_B timetable updates
# ---------------------------------------------------------------------------
$wsB = $wb.Worksheets.Item("Section_B")

$wsB.Range("B2").Value = "CS264"
$wsB.Range("E2").Value = "MA261"
$wsB.Range("F2").Value = "Free"

$wsB.Range("B3").Value = "CS264 (Tutorial)"
$wsB.Range("D3").Value = "Free"
$wsB.Range("E3").Value = "CS263"

$wsB.Range("C5").Value = "CS264"
$wsB.Range("D5").Value = "Free"
$wsB.Range("E5").Value = "CS261"
$wsB.Range("F5").Value = "Free"

$wsB.Range("C6").Value = "MA261"
$wsB.Range("D6").Value = "CS263"
$wsB.Range("F6").Value = "CS264"

$wsB.Range("B7").Value = "CS261"
$wsB.Range("C7").Value = "CS261"
$wsB.Range("E7").Value = "Free"
$wsB.Range("F7").Value = "CS263"

# ---------------------------------------------------------------------------
# Course_Summary updates (course codes renamed / renumbered + new instructors)
# ---------------------------------------------------------------------------
$wsC = $wb.Worksheets.Item("Course_Summary")

$wsC.Range("A2").Value = "MA261"
$wsC.Range("B2").Value = "Differential Equations"
$wsC.Range("F2").Value = "Dr. Anand Barangi"

$wsC.Range("A3").Value = "CS261"
$wsC.Range("B3").Value = "Operating System"
$wsC.Range("F3").Value = "Dr. Somes"

$wsC.Range("A4").Value = "CS263"
$wsC.Range("B4").Value = "Design & Analysis of Algorithms"
$wsC.Range("F4").Value = "Dr. Prabhu Prasad"

$wsC.Range("A5").Value = "CS264"
$wsC.Range("B5").Value = "Computer Networks"
$wsC.Range("F5").Value = "Dr. Prabhu Prasad"

# ---------------------------------------------------------------------------
# New sheet: Elective_Coordination (empty placeholder sheet, added at the end)
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsNew = $wb.Worksheets.Add($null, $lastSheet)
$wsNew.Name = "Elective_Coordination"
